$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title (volume/number) and date-range strings ---
$ws.Range("A8").Value = "Volume 32   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# --- Plain numeric cell updates ---
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 100
$ws.Range("N15").Value = -46.666666666666
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 37
$ws.Range("K16").Value = -5.128205128205
$ws.Range("L16").Value = -9.756097560975
$ws.Range("M16").Value = -31.481481481481
$ws.Range("N16").Value = -88.064516129032
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -18.181818181818
$ws.Range("I17").Value = 86
$ws.Range("J17").Value = 65
$ws.Range("K17").Value = 32.307692307692
$ws.Range("L17").Value = 7.5
$ws.Range("M17").Value = 48.275862068965
$ws.Range("N17").Value = -41.891891891891
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 200
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = -2.127659574468
$ws.Range("L18").Value = -22.033898305084
$ws.Range("M18").Value = -63.779527559055
$ws.Range("N18").Value = -91.943957968476
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 38.461538461538
$ws.Range("I19").Value = 188
$ws.Range("J19").Value = 220
$ws.Range("K19").Value = -14.545454545454
$ws.Range("L19").Value = -20.338983050847
$ws.Range("M19").Value = 16.049382716049
$ws.Range("N19").Value = -25.396825396825
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -14.285714285714
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = -56.25
$ws.Range("I20").Value = 65
$ws.Range("J20").Value = 110
$ws.Range("K20").Value = -40.909090909090
$ws.Range("L20").Value = 6.557377049180
$ws.Range("M20").Value = -24.418604651162
$ws.Range("N20").Value = -93.885230479774
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -7.594936708860
$ws.Range("I21").Value = 430
$ws.Range("J21").Value = 486
$ws.Range("K21").Value = -11.522633744856
$ws.Range("L21").Value = -11.340206185567
$ws.Range("M21").Value = -12.955465587044
$ws.Range("N21").Value = -81.795088907705
$ws.Range("M22").Value = -20
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 31.578947368421
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -29.523809523809
$ws.Range("I24").Value = 566
$ws.Range("J24").Value = 729
$ws.Range("K24").Value = -22.359396433470
$ws.Range("L24").Value = -30.123456790123
$ws.Range("M24").Value = -8.116883116883
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -44.444444444444
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 72
$ws.Range("H25").Value = -62.5
$ws.Range("I25").Value = 269
$ws.Range("J25").Value = 470
$ws.Range("K25").Value = -42.765957446808
$ws.Range("L25").Value = -34.549878345498
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -37.5
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 215
$ws.Range("J26").Value = 194
$ws.Range("K26").Value = 10.824742268041
$ws.Range("L26").Value = 5.392156862745
$ws.Range("M26").Value = 11.398963730569
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 100
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 20
$ws.Range("K28").Value = -23.076923076923
$ws.Range("L28").Value = 33.333333333333
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 0

# --- Cells becoming the literal text marker '***.*' (non-numeric, safe as plain .Value) ---
$ws.Range("E16").Value = "***.*"
$ws.Range("E31").Value = "***.*"

# --- Cells becoming the literal text marker '0' (ambiguous with numeric 0; force text via formula+paste-values trick to keep it a shared string without altering cell style) ---
$ws.Range("D16").Formula = '="0"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("C18").Formula = '="0"'
$ws.Range("C18").Copy()
$ws.Range("C18").PasteSpecial(-4163)
$ws.Range("D31").Formula = '="0"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
